# Refresh the "Coin"/"Link"/"Price"/"Volume(1h)" columns of the crypto
# ranking sheet with the latest scrape. A handful of rows in the 7-19 and
# 41-43 ranges also changed order (coins re-ranked by the scraper), so both
# the name/link (B/C) and the price/volume (D/E) are rewritten for those rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> column letter -> new cell text
$updates = [ordered]@{
    2 = @{ "D" = "'246.57"; "E" = "'0.97%" }
    3 = @{ "D" = "'29.42"; "E" = "'7.52%" }
    4 = @{ "D" = "'5.200"; "E" = "'3.04%" }
    5 = @{ "D" = "'0.05711"; "E" = "'0.60%" }
    6 = @{ "D" = "'6.575"; "E" = "'1.54%" }
    7 = @{ "B" = "'MXToken"; "C" = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; "D" = "'0.8582"; "E" = "'4.39%" }
    8 = @{ "B" = "'FTXToken"; "C" = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; "D" = "'0.8775"; "E" = "'4.64%" }
    9 = @{ "B" = "'WazirX"; "C" = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D" = "'0.1367"; "E" = "'3.18%" }
    10 = @{ "B" = "'MandalaExchangeToken"; "C" = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D" = "'0.07085"; "E" = "'2.32%" }
    11 = @{ "B" = "'BitrueCoin"; "C" = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; "D" = "'0.02872"; "E" = "'0.36%" }
    12 = @{ "B" = "'BitMartToken"; "C" = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; "D" = "'0.09385"; "E" = "'0.00%" }
    13 = @{ "B" = "'BitForexToken"; "C" = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; "D" = "'0.001544"; "E" = "'1.65%" }
    14 = @{ "B" = "'CoinExToken"; "C" = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; "D" = "'0.04159"; "E" = "'0.61%" }
    15 = @{ "B" = "'One"; "C" = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; "D" = "'0.0006008"; "E" = "'0.30%" }
    16 = @{ "B" = "'TigerCash"; "C" = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"; "D" = "'0.006105"; "E" = "'-0.41%" }
    17 = @{ "B" = "'UpBots"; "C" = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"; "D" = "'0.007491"; "E" = "'5,107.41%" }
    18 = @{ "B" = "'LEO"; "C" = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; "D" = "'3.485"; "E" = "'-0.67%" }
    19 = @{ "B" = "'GateToken"; "C" = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; "D" = "'3.061"; "E" = "'2.08%" }
    20 = @{ "D" = "'2.182"; "E" = "'-1.86%" }
    21 = @{ "D" = "'0.3177"; "E" = "'2.07%" }
    22 = @{ "D" = "'0.03303"; "E" = "'4.53%" }
    23 = @{ "D" = "'0.1301"; "E" = "'3.71%" }
    24 = @{ "D" = "'3.467"; "E" = "'-3.27%" }
    25 = @{ "E" = "'0.41%" }
    26 = @{ "D" = "'0.005054"; "E" = "'30.62%" }
    27 = @{ "E" = "'-0.22%" }
    28 = @{ "D" = "'0.0001209"; "E" = "'23.34%" }
    40 = @{ "D" = "'0.03746"; "E" = "'1.44%" }
    41 = @{ "B" = "'BKEXToken"; "C" = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; "D" = "'0.1074"; "E" = "'1.85%" }
    42 = @{ "B" = "'CEJI"; "C" = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; "D" = "'0.002538"; "E" = "'10.69%" }
    43 = @{ "B" = "'KickToken"; "C" = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; "D" = "'0.003469"; "E" = "'-43.62%" }
    44 = @{ "E" = "'-2.98%" }
    45 = @{ "D" = "'0.00005099"; "E" = "'-3.32%" }
    46 = @{ "E" = "'-0.03%" }
    47 = @{ "D" = "'0.07098"; "E" = "'-30.07%" }
    48 = @{ "D" = "'0.002594"; "E" = "'1.24%" }
    49 = @{ "D" = "'0.00002099"; "E" = "'-0.03%" }
    50 = @{ "D" = "'0.0001999"; "E" = "'-0.03%" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellRef = "$col$row"
        # Leading apostrophe forces text storage so "5.200" stays "5.200"
        # instead of being parsed into the float 5.2.
        $ws.Range($cellRef).Value = $updates[$row][$col]
        # Drop the quote-prefix style the text assignment implicitly applied
        # so the cell formatting matches the untouched cells around it.
        $ws.Range($cellRef).Style = "Normal"
    }
}
